$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.873.20'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.080.91'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.43'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.42'
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.079.55'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.92'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.92'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.583.93'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.708.26'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.077.75'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '486.00'
$ws.Range('E20').Value = '  +3.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.60'
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.712'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.63'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.47'
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.06'
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.95'
$ws.Range('E26').Value = '  -2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.77'
$ws.Range('E27').Value = '  +8.86%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.51'
$ws.Range('E29').Value = '  +2.66%  '
$ws.Range('E30').Value = '  +2.46%  '
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.38'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('E35').Value = '  +1.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0826'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.09'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.26'
$ws.Range('E39').Value = '  -3.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.30'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.78'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '440.57'
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.292'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.113'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0367'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.841.91'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.56'
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.96'
$ws.Range('E48').Value = '  +2.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.63'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.26'
$ws.Range('E51').Value = '  +0.18%  '
